$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# ------------------------------------------------------------------
# The sheet currently has:
#   row 4       : last data row (template row we will clone for the new items)
#   row 5       : totals row              (K5:N5 merged)
#   row 6       : footer row              (A6:E6, F6:G6, I6:N6 merged)
# Target layout adds 5 new data rows (new row 5..9) and pushes the totals
# row to row 10 and the footer row to row 11.
#
# We relocate rows 5/6 -> 10/11 manually (copy values + formats, rebuild
# the merges) instead of using EntireRow.Insert(), which would otherwise
# make Excel auto-generate extra unused style records for the freshly
# inserted blank rows. We also merge the *destination* range BEFORE
# pasting the formats onto it, so every cell in the merged range ends up
# sharing the same style index (matching the source layout), instead of
# the style getting fragmented per-cell by merging after formatting.
# ------------------------------------------------------------------

# Remember old row 5 / row 6 content before we touch anything.
$totalValue = $ws.Range("K5").Value()
$footerA = $ws.Range("A6").Value()
$footerF = $ws.Range("F6").Value()
$footerI = $ws.Range("I6").Value()

# Break the old merges so we can freely copy/move the ranges.
$ws.Range("K5:N5").UnMerge()
$ws.Range("A6:E6").UnMerge()
$ws.Range("F6:G6").UnMerge()
$ws.Range("I6:N6").UnMerge()

# --- Move totals row (old row 5) down to row 10 ---
$ws.Range("K10:N10").Merge()
$ws.Range("K5:N5").Copy()
$ws.Range("K10:N10").PasteSpecial($xlPasteFormats)
$ws.Range("K10").Value = $totalValue
$ws.Rows.Item(10).RowHeight = 26.25

# --- Move footer row (old row 6) down to row 11 ---
$ws.Range("A11:E11").Merge()
$ws.Range("A6:E6").Copy()
$ws.Range("A11:E11").PasteSpecial($xlPasteFormats)
$ws.Range("A11").Value = $footerA

$ws.Range("F11:G11").Merge()
$ws.Range("F6:G6").Copy()
$ws.Range("F11:G11").PasteSpecial($xlPasteFormats)
$ws.Range("F11").Value = $footerF

$ws.Range("H6").Copy()
$ws.Range("H11").PasteSpecial($xlPasteFormats)

$ws.Range("I11:N11").Merge()
$ws.Range("I6:N6").Copy()
$ws.Range("I11:N11").PasteSpecial($xlPasteFormats)
$ws.Range("I11").Value = $footerI

$ws.Rows.Item(11).RowHeight = 16.5

# Clear the now-stale old row 5 / row 6 cells (they become part of the new
# item rows below and will be repopulated from the row-4 template).
$ws.Range("A5:N6").Clear()

# --- New row data: row number -> [A(seq), B(name), H(balance), L(price), N(count), row height] ---
$rows = @(
    @{ Row = 5; A = 2; B = "DANSET 8MG/4ML 3 AMP."; H = "0:1";  L = 95; N = "0:3"; Height = 25.5 },
    @{ Row = 6; A = 3; B = "جهاز محلول ";            H = "10:0"; L = 20; N = "1:0"; Height = 24.75 },
    @{ Row = 7; A = 4; B = "سرنجات 5 سم";            H = "-1:0"; L = 2;  N = "1:0"; Height = 25.5 },
    @{ Row = 8; A = 5; B = "كالونا ";                H = "-1:0"; L = 15; N = "1:0"; Height = 25.5 },
    @{ Row = 9; A = 6; B = "محلول ملح";              H = "27:0"; L = 48; N = "2:0"; Height = 24.75 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Merge the destination ranges first so the subsequent format-paste
    # lands a single consistent style across every cell in the group
    # (matching row 4's column-style groups: A | B:G | H:K | L:M | N),
    # without creating new style entries in the workbook.
    $ws.Range("B$rowNum`:G$rowNum").Merge()
    $ws.Range("H$rowNum`:K$rowNum").Merge()
    $ws.Range("L$rowNum`:M$rowNum").Merge()

    $ws.Range("A4").Copy()
    $ws.Range("A$rowNum").PasteSpecial($xlPasteFormats)

    $ws.Range("B4:G4").Copy()
    $ws.Range("B$rowNum`:G$rowNum").PasteSpecial($xlPasteFormats)

    $ws.Range("H4:K4").Copy()
    $ws.Range("H$rowNum`:K$rowNum").PasteSpecial($xlPasteFormats)

    $ws.Range("L4:M4").Copy()
    $ws.Range("L$rowNum`:M$rowNum").PasteSpecial($xlPasteFormats)

    $ws.Range("N4").Copy()
    $ws.Range("N$rowNum").PasteSpecial($xlPasteFormats)

    # Values
    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("H$rowNum").Value = $r.H
    $ws.Range("L$rowNum").Value = $r.L
    $ws.Range("N$rowNum").Value = $r.N

    # Row height
    $ws.Rows.Item($rowNum).RowHeight = $r.Height
}

# --- Update the totals row (row 10): sum of all L column prices ---
$ws.Range("K10").Value = 199

Write-Host "edit complete"
